# New crime data collected - weekly refresh of the 123rd Precinct CompStat
# report: header volume/date text, and the weekly / 28-day / YTD / historical
# crime-complaint figures (and their derived %-change columns) for rows 16-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text (rich-text shared strings): issue number and the week-covered
# date range.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# ---------------------------------------------------------------------------
# Row 16 - Rape
# ---------------------------------------------------------------------------
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 15
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -6.25
$ws.Range("L16").Value = 36.363636363636
$ws.Range("M16").Value = 15.384615384615
$ws.Range("N16").Value = -42.307692307692

# ---------------------------------------------------------------------------
# Row 17 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 9
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 57
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = 9.615384615384
$ws.Range("L17").Value = 14
$ws.Range("M17").Value = 159.090909090909
$ws.Range("N17").Value = -5

# ---------------------------------------------------------------------------
# Row 18 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("L18").Value = -3.125
$ws.Range("M18").Value = -49.180327868852
$ws.Range("N18").Value = -80

# ---------------------------------------------------------------------------
# Row 19 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -19.354838709677
$ws.Range("I19").Value = 141
$ws.Range("J19").Value = 169
$ws.Range("K19").Value = -16.568047337278
$ws.Range("L19").Value = -10.191082802547
$ws.Range("M19").Value = 98.591549295774
$ws.Range("N19").Value = 30.555555555555

# ---------------------------------------------------------------------------
# Row 20 - Gr. Larceny (C20 goes from a numeric 2 to the text "0" marker,
# copying D20's formatting/value so the style + shared-string text match).
# ---------------------------------------------------------------------------
$ws.Range("D20").Copy($ws.Range("C20"))
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("L20").Value = -65.909090909090
$ws.Range("N20").Value = -95.901639344262

# ---------------------------------------------------------------------------
# Row 21 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 11.111111111111
$ws.Range("F21").Value = 45
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = 2.272727272727
$ws.Range("I21").Value = 263
$ws.Range("J21").Value = 280
$ws.Range("K21").Value = -6.071428571428
$ws.Range("L21").Value = -11.148648648648
$ws.Range("M21").Value = 38.421052631578
$ws.Range("N21").Value = -63.421418636995

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 366.666666666667
$ws.Range("F24").Value = 36
$ws.Range("G24").Value = 24
$ws.Range("H24").Value = 50
$ws.Range("I24").Value = 239
$ws.Range("J24").Value = 208
$ws.Range("K24").Value = 14.903846153846
$ws.Range("L24").Value = -8.076923076923
$ws.Range("M24").Value = -11.152416356877

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("E25").Value = 800
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 144.444444444444
$ws.Range("I25").Value = 140
$ws.Range("J25").Value = 97
$ws.Range("K25").Value = 44.329896907216
$ws.Range("L25").Value = 21.739130434782

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 113
$ws.Range("J26").Value = 89
$ws.Range("K26").Value = 26.966292134831
$ws.Range("L26").Value = 8.653846153846
$ws.Range("M26").Value = 0.892857142857

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50

# ---------------------------------------------------------------------------
# Rows 29/30 - Shooting Vic./Inc. (L29/L30 go from the "***.*" text marker
# to an actual -100 number, copying N29/N30's formatting so the style
# matches the numeric format used elsewhere in the row).
# ---------------------------------------------------------------------------
$ws.Range("N29").Copy($ws.Range("L29"))
$ws.Range("N30").Copy($ws.Range("L30"))

# ---------------------------------------------------------------------------
# Column widths for columns E and H (5 and 8) shrink to match the other
# %-chg columns (F, G, I, J).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
